# Adapt column header formatting to respective input file names (#7):
# header suffix "_old" -> "_FV2404", "_new" -> "_FV2410", then turn the
# sheet's data range into an Excel Table and freeze the header row so the
# column names stay visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- Rename header cells in row 1 -----------------------------------------
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2404")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2410")
    }
}

# --- Freeze the header row -------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Convert the data range into an Excel Table -----------------------------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"
